# Generate Report for handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - New "Latest Target File" / "Latest Handback File" hyperlink columns (E/F) are
#    populated for the two tracked source files on both the zh-cn and de-de sheets
#  - The "Latest Handback DateTime" column (G) is updated with the real handback time

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet: Status columns for zh-cn / de-de mirror the same shared text ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $statusText
$wsZh.Range("B3").Value = $statusText

# Latest Target File (E) / Latest Handback File (F) hyperlinks for row 2 (6cc58f80...)
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/f91d7ef379800eb96c02efa54de5d21f674dae31/e2e/6cc58f80-6c75-4310-99b5-213a766b612f.md", "", "", "6cc58f80-6c75-4310-99b5-213a766b612f.md") | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/835d996f165ce5765eb9643a6987bc8acde3d8b4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/6cc58f80-6c75-4310-99b5-213a766b612f.da90bb23d96255d1adb05ca2c065052b393dd1b9.zh-cn.xlf", "", "", "6cc58f80-6c75-4310-99b5-213a766b612f.da90bb23d96255d1adb05ca2c065052b393dd1b9.zh-cn.xlf") | Out-Null

# Latest Target File (E) / Latest Handback File (F) hyperlinks for row 3 (72168430...)
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/f91d7ef379800eb96c02efa54de5d21f674dae31/e2e/72168430-c069-4eac-ae40-91b6ac97514e.md", "", "", "72168430-c069-4eac-ae40-91b6ac97514e.md") | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/835d996f165ce5765eb9643a6987bc8acde3d8b4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/72168430-c069-4eac-ae40-91b6ac97514e.3496b935a40f29a1337b095f59d923a35022b25f.zh-cn.xlf", "", "", "72168430-c069-4eac-ae40-91b6ac97514e.3496b935a40f29a1337b095f59d923a35022b25f.zh-cn.xlf") | Out-Null

# Latest Handback DateTime (G) for both tracked rows
$wsZh.Range("G2").Value = "2016-01-18 05:06:52"
$wsZh.Range("G3").Value = "2016-01-18 05:06:52"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $statusText
$wsDe.Range("B3").Value = $statusText

# Latest Target File (E) / Latest Handback File (F) hyperlinks for row 2 (6cc58f80...)
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/f91d7ef379800eb96c02efa54de5d21f674dae31/e2e/6cc58f80-6c75-4310-99b5-213a766b612f.md", "", "", "6cc58f80-6c75-4310-99b5-213a766b612f.md") | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d41715821140a2a79e97b1539c41e8f872597b88/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/6cc58f80-6c75-4310-99b5-213a766b612f.da90bb23d96255d1adb05ca2c065052b393dd1b9.de-de.xlf", "", "", "6cc58f80-6c75-4310-99b5-213a766b612f.da90bb23d96255d1adb05ca2c065052b393dd1b9.de-de.xlf") | Out-Null

# Latest Target File (E) / Latest Handback File (F) hyperlinks for row 3 (72168430...)
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/f91d7ef379800eb96c02efa54de5d21f674dae31/e2e/72168430-c069-4eac-ae40-91b6ac97514e.md", "", "", "72168430-c069-4eac-ae40-91b6ac97514e.md") | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d41715821140a2a79e97b1539c41e8f872597b88/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/72168430-c069-4eac-ae40-91b6ac97514e.3496b935a40f29a1337b095f59d923a35022b25f.de-de.xlf", "", "", "72168430-c069-4eac-ae40-91b6ac97514e.3496b935a40f29a1337b095f59d923a35022b25f.de-de.xlf") | Out-Null

# Latest Handback DateTime (G) for both tracked rows
$wsDe.Range("G2").Value = "2016-01-18 05:07:13"
$wsDe.Range("G3").Value = "2016-01-18 05:07:13"
